$wb = $excel.ActiveWorkbook

# The Overview sheet mirrors each language sheet's Status in row 2
# ("Ready for handoff" -> "Handoff transform failed"); keep it in sync so
# the shared string is updated in place rather than leaving a duplicate.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value2 = "Handoff transform failed"
$overview.Range("C2").Value2 = "Handoff transform failed"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status: "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value2 = "Handoff transform failed"

    # Remove the "Latest Handoff File" hyperlink + value in C2 entirely.
    $linksToRemove = @()
    foreach ($link in $ws.Hyperlinks) {
        if ($link.Range.Address() -eq '$C$2') {
            $linksToRemove += $link
        }
    }
    foreach ($link in $linksToRemove) {
        $link.Delete()
    }
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime reset to the epoch placeholder.
    $ws.Range("D2").Value2 = "0001-01-01 00:00:00"

    # Handoff Reason: "Include" -> "Ignored"
    $ws.Range("H2").Value2 = "Ignored"
}
